# Applies the "mise à jour avant publication" edit to dist_coûts.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Vertical-alignment swaps (these are the only difference between the two
#    date cell styles used on this sheet: one centers vertically, the other
#    uses the default/bottom alignment). Setting .VerticalAlignment directly
#    can corrupt the shared date NumberFormat in this runtime, so instead we
#    copy the *format only* from a donor cell that already carries the style
#    we want, which correctly reuses the existing style entries.
# ---------------------------------------------------------------------------

# D12:D17, F2:F25, H2:H25 : default/bottom -> centered (style used by D2, etc.)
$ws.Range("D2").Copy()
$ws.Range("D12:D17").PasteSpecial(-4122)
$ws.Range("F2:F25").PasteSpecial(-4122)
$ws.Range("H2:H25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Q2:Q18 : centered -> default/bottom (style used by S2, etc.)
$ws.Range("S2").Copy()
$ws.Range("Q2:Q18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Value updates
# ---------------------------------------------------------------------------

# K2:K9 date values 45352 -> 45748
$ws.Range("K2:K9").Value = 45748

# J6, J7, J8, J9 -> 9999
$ws.Range("J6").Value = 9999
$ws.Range("J7").Value = 9999
$ws.Range("J8").Value = 9999
$ws.Range("J9").Value = 9999

# Y8, Z8, AA8 -> 0
$ws.Range("Y8").Value = 0
$ws.Range("Z8").Value = 0
$ws.Range("AA8").Value = 0

# C11 -> 10.8
$ws.Range("C11").Value = 10.8

# ---------------------------------------------------------------------------
# 3. Sheet view / selection
# ---------------------------------------------------------------------------
$ws.Range("J6").Select()
